$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New coded-segment rows (236-241) appended by the latest mex run / QA pass
# over segments, drugs and bacteria.
$newRows = @(
    @{ Row=236; D="2697"; E="Event month"; F="2: 2572"; G="2: 2574"; I="May";    J=3; K=0.011559; M="11/14/18 11:31:00" },
    @{ Row=237; D="2697"; E="Event month"; F="2: 2580"; G="2: 2585"; I="August"; J=6; K=0.023118; M="11/14/18 11:31:00" },
    @{ Row=238; D="2697"; E="Event year";  F="2: 2587"; G="2: 2590"; I="2002";   J=4; K=0.015412; M="11/14/18 11:31:00" },
    @{ Row=239; D="3651"; E="Event year";  F="6: 665";  G="6: 668";  I="2013";   J=4; K=0.008877; M="11/14/18 11:33:00" },
    @{ Row=240; D="3651"; E="Event year";  F="6: 655";  G="6: 658";  I="2012";   J=4; K=0.008877; M="11/14/18 11:33:00" },
    @{ Row=241; D="3910"; E="Event year";  F="4: 2261"; G="4: 2264"; I="2006";   J=4; K=0.020517; M="11/14/18 11:35:00" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $prev = $rowNum - 1

    # Clone formatting (styles + row height) from the row directly above.
    $src = $ws.Range("A" + $prev + ":M" + $prev)
    $dst = $ws.Range("A" + $rowNum + ":M" + $rowNum)
    $src.Copy($dst)
    $ws.Rows($rowNum).RowHeight = 16

    $ws.Cells.Item($rowNum, 1).Value = "●"
    $ws.Cells.Item($rowNum, 2).Value = ""
    $ws.Cells.Item($rowNum, 3).Value = ""

    # Columns D and I occasionally hold purely-numeric text (document/segment
    # ids, years); force them to stay text with a leading apostrophe and then
    # re-paste the (now-blank, still clean) column-B formatting so the style
    # doesn't pick up a stray quote-prefix format.
    $ws.Cells.Item($rowNum, 4).Formula = "'" + $r.D
    $ws.Range("B" + $rowNum).Copy()
    $ws.Range("D" + $rowNum).PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = 0

    $ws.Cells.Item($rowNum, 9).Formula = "'" + $r.I
    $ws.Range("B" + $rowNum).Copy()
    $ws.Range("I" + $rowNum).PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = "Sonia"
    $ws.Cells.Item($rowNum, 13).Value = $r.M
}
